# "Common: Another set of data for import"
# Adds a new vendor ("Rocket Girl") to the "vendors" sheet, keeping the
# existing alphabetically-sorted list in order (it sorts between
# "Ripe Vapes" and "Samsung", landing on row 63 and pushing the rest down).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vendors")

# Insert a new row at the correct sorted position and fill in the value.
$ws.Rows("63:63").Insert()
$ws.Range("A63").Value = "Rocket Girl"

# Re-apply the sort over the (now one-row-bigger) data range so the sheet's
# remembered sort state follows the new extent.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A90:A96")) | Out-Null
$ws.Sort.SetRange($ws.Range("A2:A96"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Match the author's final scroll position/selection on the sheet.
$ws.Range("A81").Select()
